$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = "-"

# Row 3
$ws.Range("B3").Value = "[Carlos-Tornearia, Victor-Ajustagem, -, -]"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Gilberto-Mec. Tec. Res. M"
$ws.Range("E3").Value = "[Elcio Dec.-Des. Maq. Cad._T2, -]"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "[Carlos-Tornearia, Victor-Ajustagem, Emerson-Elet. Dig. Bas., Elaine-Metalografia]"
$ws.Range("C4").Value = "[-, Elcio Dec.-Des. Maq. Cad._T1]"
$ws.Range("E4").Value = "[Elcio Dec.-Des. Maq. Cad._T2, -]"
$ws.Range("F4").Value = "[Carlos-Tornearia, Emerson-Elet. Dig. Bas., Elaine-Metalografia, Victor-Ajustagem]"

# Row 6
$ws.Range("B6").Value = "[-, Victor-Ajustagem, Emerson-Elet. Dig. Bas., Elaine-Metalografia]"
$ws.Range("C6").Value = "[-, Elcio Dec.-Des. Maq. Cad._T1]"
$ws.Range("D6").Value = "Gilberto-Mec. Tec. Res. M"
$ws.Range("E6").Value = "[Elcio Dec.-Des. Maq. Cad._T2, -]"
$ws.Range("F6").Value = "[Carlos-Tornearia, Emerson-Elet. Dig. Bas., Elaine-Metalografia, -]"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "[-, Elcio Dec.-Des. Maq. Cad._T1]"
$ws.Range("D7").Value = "Maria Celeste-Máquinas Térmicas e de Fl"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("D8").Value = "Maria Celeste-Máquinas Térmicas e de Fl"
$ws.Range("F8").Value = "-"
